# Updated symbol list on Wed Jan 18 07:58:46 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) quotes for the
# coin rows on the active sheet to the latest scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Values must be written as TEXT
# (the sheet stores Price/Volume as strings, e.g. "301.66" / "0.66%"),
# so each cell is temporarily forced to Text format, written, then the
# number-format override is cleared again (leaving no style change behind).
$updates = [ordered]@{
    "D2"  = "301.66";  "E2"  = "0.66%"
    "D3"  = "32.76";   "E3"  = "4.16%"
    "D4"  = "4.943";   "E4"  = "-2.76%"
    "D5"  = "0.07756"; "E5"  = "-1.28%"
    "D6"  = "1.982";   "E6"  = "-12.58%"
    "D7"  = "7.848";   "E7"  = "0.44%"
    "D8"  = "3.797";   "E8"  = "-1.31%"
    "D9"  = "0.9209";  "E9"  = "0.11%"
    "D10" = "0.1770"
    "D11" = "0.07879"; "E11" = "3.76%"
    "D12" = "0.08572"; "E12" = "-6.71%"
    "D13" = "0.03142"
    "E14" = "0.03%"
    "D15" = "0.001510"; "E15" = "0.41%"
    "D16" = "0.005874"; "E16" = "-0.40%"
    "D18" = "2.153";    "E18" = "-4.01%"
    "D19" = "0.3338";   "E19" = "2.02%"
    "E20" = "-0.02%"
    "D21" = "4.306";    "E21" = "10.11%"
    "D22" = "0.1991";   "E22" = "16.26%"
    "D23" = "0.04559";  "E23" = "-1.25%"
    "D24" = "0.001225"; "E24" = "-2.22%"
    "D25" = "0.004432"; "E25" = "-0.82%"
    "D26" = "0.0001250"; "E26" = "0.10%"
    "D39" = "0.01716";  "E39" = "-1.17%"
    "D40" = "0.04720";  "E40" = "2.14%"
    "D41" = "0.008018"; "E41" = "13.61%"
    "E42" = "-0.26%"
    "D43" = "0.002341"; "E43" = "6.96%"
    "D44" = "0.01045";  "E44" = "7.31%"
    "D45" = "0.00006237"; "E45" = "-0.47%"
    "D46" = "0.00000000750"; "E46" = "0.12%"
    "D47" = "0.8234";   "E47" = "10.35%"
    "D48" = "0.003101"; "E48" = "-61.17%"
    "D49" = "0.00002101"; "E49" = "0.12%"
    "D50" = "0.0002001"; "E50" = "0.12%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
